$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2 through 76, replacing the
# previous Strike# derived values with the regenerated K values.
$gValues = @(2,2,0,3,1,1,1,1,1,0,1,1,0,0,0,2,1,2,0,3,1,0,0,1,1,0,1,2,0,1,1,1,1,0,1,2,1,0,3,3,1,1,1,3,1,2,1,2,1,2,5,2,2,1,3,5,2,2,3,1,0,1,1,2,0,1,0,0,0,0,1,3,0,2,0)

$startRow = 2
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
}
